$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2025 UiPath AI and Automation Trends"
$ws.Range("A2").Value = "Unlocking data potential:"
$ws.Range("A3").Value = "Testing in the age of generative AI"
$ws.Range("A4").Value = "Process mining and the path to digital transformation"
$ws.Range("A5").Value = "Turn AI potential into AI results: 8 steps to success in banking"
$ws.Range("A6").Value = "5 Ways UiPath Test Suite can improve your testing"
$ws.Range("A9").Value = "Turn AI potential into AI results: 8 steps to success in manufacturing"
$ws.Range("A10").Value = "Intelligent document processing: Build your own, or buy as a service?"
$ws.Range("A11").Value = "Turn AI potential into AI results: 8 steps to success for healthcare payers"
$ws.Range("A12").Value = "The Ultimate Guide to Intelligent Document Processing"
